# Apply edits described by the commit diff: update existing rows 2-21,
# append two new rows (22, 23) and move the "Total" row to row 24 with
# recomputed totals. All values are textual (Brazilian-formatted numbers
# such as "29.271" must stay strings, not be parsed as numeric), so each
# cell is forced to text format before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2..24, columns A..P. $null means "leave the cell empty".
$rowsData = @(
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA SJ DO RIO PRETO', '0081/23C', 'VIBRA - STI VILLE  - BL 06 - 0,8165', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '29.271', '24.002', '0', '0', '0', '0', '29.271', '24.002', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA BARUERI', '0102/23B', 'VIBRA - WECO MADELEINE - BL 2 - 0,8304', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '39', '35', '0', '0', '0', '0', '39', '35', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA BARUERI', '0278/23F', 'TRANSF VIBRA.CUBATÃO X V.BARUERI CONF NF 1273775-0 29.08.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '900', '3.489', '0', '0', '0', '0', '900', '3.489', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA BARUERI', '0283/23C', 'VIBRA X VIBRA BARUERI CONF NF 1278522 20.09.23', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '3.240', '2.430', '0', '0', '0', '0', '3.240', '2.430', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA SAO PAULO', '0283/23D', 'VIBRA CUB X V. SP CONF NF 1279755-0 25.09.23', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '6.677', '5.371', '0', '0', '0', '0', '6.677', '5.371', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA BARUERI', '0283/23E', 'VIBRA SP X V. BARUERI  CONF NF 3103142-1 26.09.23', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '9.938', '8.170', '0', '0', '0', '0', '9.938', '8.170', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA CUBATAO', '0283/23G', 'TRANSF VIBRA PAU X  VIBRA CUB CONF NF 3749230-0 04.10.23', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '86', '4.190', '0', '0', '0', '0', '86', '4.190', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA RIBEIRAO PRETO', '0283/23H', 'VIBRA - AQUARIUS - NF 1272825 DENSIDADE 0,8290', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '16.640', '13.385', '0', '0', '0', '0', '16.640', '13.385', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA CUBATAO', '0283/23I', 'TRANSF VIBRA PAULINIA CONF NF 3749795-1 05.09.23', 'OLEO DIESEL S500', 'ONU 1202, OLEO DIESEL S500', '12.421', '6.166', '0', '0', '0', '0', '12.421', '6.166', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA RIBEIRAO PRETO', '0307/23J', 'TRANSF VIBRA CUB X V.RIB PRETO CONF NF 1281710 - 0 03.10.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '8.967', '8.257', '0', '0', '0', '0', '8.967', '8.257', '0')
    ,('QUEBRA', 'MERCADO INTERNO', 'VIBRA SAO PAULO', '0307/23M', 'TRANSF  VIBRA PAU CONF NF 3750527-1 06.10.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '34.074', '28.080', '0', '0', '0', '0', '34.074', '28.080', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PAULINIA', '0174/23Z', 'TRANSF VIBRA CUBATÃO CONF NFE 1264138-0 DE 17.07.2023', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '73', '62', '0', '0', '0', '0', '73', '62', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PRESIDENTE PRUDENTE', '0259/23G', 'TRANSF VIBRA CUB X PRE.PRUDENTE CONF NF 1276885-0 12.09.2023', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '13.503', '11.526', '0', '0', '0', '0', '13.503', '11.526', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PAULINIA', '0307/23L', 'TRANSF  VIBRA CUB CONF NF1282508 1 06.10.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '45.628', '38.507', '0', '0', '0', '0', '45.628', '38.507', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA SAO PAULO', '0307/23M', 'TRANSF  VIBRA PAU CONF NF 3750527-1 06.10.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '72.525', '60.000', '0', '0', '0', '0', '72.525', '60.000', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PAULINIA', '0307/23N', 'AJUSTE DE SALDO (CONF NF 3759087-1 24.10.23)', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '8.052', '6.764', '0', '0', '0', '0', '8.052', '6.764', '0')
    ,('TQ-03', 'IMPORTACAO COMUM', 'VIBRA CUBATAO', '0322/23A', 'TORM ALTLANTIC - VIBRA CUBATÃO - BL 02 DENS:0,8306', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '21.161', '3.597', '0', '0', '0', '0', '21.161', '3.597', '0')
    ,('TQ-03', 'IMPORTACAO COMUM', 'VIBRA CUBATAO', '0322/23B', 'TORM ALTLANTIC - VIBRA CUBATÃO - BL 08 DENS: FALTA', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '0', '0', '0', '5.006.360', '4.158.176', '0', '5.006.360', '4.158.176', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PAULINIA', '0322/23C', 'TRANSF VIBRA CUB X PAULINIA CONF NF 1286527-0 25.10.23', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '186.385', '154.810', '0', '0', '0', '0', '186.385', '154.810', '0')
    ,('TQ-03', 'MERCADO INTERNO', 'VIBRA PRESIDENTE PRUDENTE', '0322/23D', 'TORM ALTLANTIC - VIBRA CUBATÃO - BL 02 DENS:0,8306', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '5.000.000', '4.166.500', '0', '0', '0', '0', '5.000.000', '4.166.500', '0')
    ,('TQ-07', 'IMPORTACAO COMUM', 'VIBRA CUBATAO', '0322/23', 'TORM ALTLANTIC - VIBRA CUBATÃO - BL 01 DENS: 0,8306', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '722', '598', '0', '0', '0', '0', '722', '598', '0')
    ,('TQ-07', 'IMPORTACAO COMUM', 'VIBRA CUBATAO', '0322/23A', 'TORM ALTLANTIC - VIBRA CUBATÃO - BL 02 DENS:0,8306', 'OLEO DIESEL S10', 'ONU 1202, DIESEL S10', '4.191.190', '3.481.102', '0', '0', '0', '0', '4.191.190', '3.481.102', '0')
    ,('Z', $null, $null, 'Total', $null, $null, $null, '9.661.492', '8.027.041', '0', '5.006.360', '4.158.176', '0', '14.667.852', '12.185.217', '0')
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rowsData[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $val = $rowValues[$c - 1]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($r, $c)
            $cell.NumberFormat = "@"
            $cell.Value = $val
        }
    }
}

